$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 0.000625
$ws.Range("K2").Value = 4265
$ws.Range("L2").Value = 0.008529999999999999
